# edit.ps1 - applies the diff changes to Диплом.pptx
#
# Summary of changes:
#  1. Slide 5, shape "TextBox 2": grow the shape height (cy 6125523 -> 6186309)
#     and split the run "Заведующий хозяйством (далее завхоз) относится к
#     категории материально-ответственных служащих" into four runs that drop
#     the word "далее ": "Заведующий " / "хозяйством " / "(завхоз" /
#     ") относится к категории материально-ответственных служащих".
#  2. Slide 7, title shape: "1С: Предприятие" -> "1С:Предприятие" (drop the
#     space after the colon).
#  3. Slide 7, body shape: split "Для разработки программного продукта была
#     выбрана платформа 1С Предприятие." into three runs so that
#     "1С Предприятие" becomes "1С:Предприятие" (insert the colon).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Slide 5 - "TextBox 2" shape
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(2)

# Grow the shape's height from 6125523 EMU to 6186309 EMU (482.324... pt ->
# 487.111 pt); width/position are untouched.
$shp5.Height = 487.111

$tf5 = $shp5.TextFrame2
$tr5 = $tf5.TextRange
$para5 = $tr5.Paragraphs(2)
$run5 = $para5.Runs(1)
$runStart = $run5.Start

$newText5 = "Заведующий хозяйством (завхоз) относится к категории материально-ответственных служащих"
$run5.Text = $newText5

# Force the single run to split into four runs at the word boundaries below
# (touching Font.Size with its current value is enough to break the run
# without altering the visible formatting).
$c1 = $tr5.Characters($runStart, 11)          # "Заведующий "
$c1.Font.Size = 22
$c2 = $tr5.Characters($runStart + 11, 11)     # "хозяйством "
$c2.Font.Size = 22
$c3 = $tr5.Characters($runStart + 22, 7)      # "(завхоз"
$c3.Font.Size = 22

# ---------------------------------------------------------------------
# 2. Slide 7 - title shape ("1С: Предприятие" -> "1С:Предприятие")
# ---------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$shpTitle = $s7.Shapes.Item(1)
$tfTitle = $shpTitle.TextFrame2
$trTitle = $tfTitle.TextRange
$paraTitle = $trTitle.Paragraphs(1)
$runTitle = $paraTitle.Runs(1)
$runTitle.Text = "1С:Предприятие"

# ---------------------------------------------------------------------
# 3. Slide 7 - body shape ("1С Предприятие" -> "1С:Предприятие")
# ---------------------------------------------------------------------
$shpBody = $s7.Shapes.Item(2)
$tfBody = $shpBody.TextFrame2
$trBody = $tfBody.TextRange
$paraBody = $trBody.Paragraphs(1)
$runBody = $paraBody.Runs(1)
$bodyStart = $runBody.Start

$newText7 = "Для разработки программного продукта была выбрана платформа 1С:Предприятие."
$runBody.Text = $newText7

$b1 = $trBody.Characters($bodyStart, 60)          # "Для ... платформа "
$b1.Font.Size = 22
$b2 = $trBody.Characters($bodyStart + 60, 14)     # "1С:Предприятие"
$b2.Font.Size = 22
